# Updated cryptos list on Sat Apr 13 03:17:04 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "66.513.55"
$ws.Range("E2").Value = "  -5.95%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.204.36"
$ws.Range("E3").Value = "  -9.07%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'575.38"
$ws.Range("E5").Value = "  -5.74%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'148.87"
$ws.Range("E6").Value = "  -14.34%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.196.52"
$ws.Range("E8").Value = "  -9.16%  "

# Row 9 - XRP
$ws.Range("D9").Value = "'0.538"
$ws.Range("E9").Value = "  -11.67%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -13.63%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'6.52"
$ws.Range("E11").Value = "  -11.50%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.493"
$ws.Range("E12").Value = "  -16.17%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "'37.93"
$ws.Range("E13").Value = "  -18.48%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -12.78%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.724.48"
$ws.Range("E15").Value = "  -9.04%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "66.506.08"
$ws.Range("E16").Value = "  -5.96%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.210.50"
$ws.Range("E17").Value = "  -8.85%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "'536.68"
$ws.Range("E18").Value = "  -12.77%  "

# Row 19 - TRON
$ws.Range("D19").Value = "'0.113"
$ws.Range("E19").Value = "  -6.59%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "'7.06"
$ws.Range("E20").Value = "  -16.37%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'14.94"
$ws.Range("E21").Value = "  -16.01%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "'0.749"
$ws.Range("E22").Value = "  -15.38%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "'7.64"
$ws.Range("E23").Value = "  -14.94%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'84.90"
$ws.Range("E24").Value = "  -13.41%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("E25").Value = "  -16.22%  "

# Row 26 - Dai
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27 - PancakeSwap
$ws.Range("D27").Value = "'3.10"
$ws.Range("E27").Value = "  -17.76%  "

# Row 28 (ranking reshuffle)
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'7.94"
$ws.Range("E28").Value = "  -13.21%  "

# Row 29 (ranking reshuffle)
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'28.91"
$ws.Range("E29").Value = "  -14.24%  "

# Row 30 (ranking reshuffle)
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'2.11"
$ws.Range("E30").Value = "  -18.91%  "

# Row 31 - Stacks
$ws.Range("E31").Value = "  -16.59%  "

# Row 32 - Mantle
$ws.Range("E32").Value = "  -14.07%  "

# Row 33 - Bittensor
$ws.Range("D33").Value = "'529.60"
$ws.Range("E33").Value = "  -14.12%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "'6.44"
$ws.Range("E34").Value = "  -21.06%  "

# Row 35 - NEARProtocol
$ws.Range("D35").Value = "'5.62"
$ws.Range("E35").Value = "  -18.06%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.15%  "

# Row 37 - OKB
$ws.Range("D37").Value = "'52.80"
$ws.Range("E37").Value = "  -7.41%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "'0.0416"
$ws.Range("E38").Value = "  -12.32%  "

# Row 39 - Hedera
$ws.Range("D39").Value = "'0.0839"
$ws.Range("E39").Value = "  -16.52%  "

# Row 40 - Cosmos
$ws.Range("D40").Value = "'8.98"
$ws.Range("E40").Value = "  -17.34%  "

# Row 41 - Kaspa
$ws.Range("D41").Value = "'0.123"
$ws.Range("E41").Value = "  -14.98%  "

# Row 42 - Maker
$ws.Range("D42").Value = "2.897.58"
$ws.Range("E42").Value = "  -14.12%  "

# Row 43 - dogwifhat
$ws.Range("D43").Value = "'2.56"
$ws.Range("E43").Value = "  -27.07%  "

# Row 44 - PEPE
$ws.Range("D44").Value = "0.0₃0579"
$ws.Range("E44").Value = "  -21.54%  "

# Row 45 - TheGraph
$ws.Range("D45").Value = "'0.257"
$ws.Range("E45").Value = "  -17.77%  "

# Row 47 - ThetaToken
$ws.Range("D47").Value = "'2.31"
$ws.Range("E47").Value = "  -21.70%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "'25.44"
$ws.Range("E48").Value = "  -21.22%  "

# Row 49 - Fetch.AI
$ws.Range("D49").Value = "'2.07"
$ws.Range("E49").Value = "  -19.69%  "

# Row 50 - Monero
$ws.Range("D50").Value = "'123.35"
$ws.Range("E50").Value = "  -7.84%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -14.35%  "

